# NIT-9009152403.xlsx: "Elimna EC anteriores y se agregan nuevos,
# se modifica base de datos"
#
# The underlying database rows for "Periodo Mora" (col E) / "Valor Mora"
# (col F) were refreshed: the period that used to be listed first (1607,
# with mora value 19305) now belongs to the last data row, and the period
# that used to be listed last (1609, with mora value 27578) now belongs to
# the first data row. The middle row (period 1608 / value 27578) is
# unchanged. Net effect: row 16 and row 18 trade their E:F contents while
# row 17 stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current "Periodo Mora" / "Valor Mora" pairs for the first and
# last data rows before overwriting anything.
$periodoRow16 = $ws.Range("E16").Value()
$valorRow16   = $ws.Range("F16").Value()
$periodoRow18 = $ws.Range("E18").Value()
$valorRow18   = $ws.Range("F18").Value()

# Swap them: row 16 gets what used to be in row 18, and vice versa.
$ws.Range("E16").Value = $periodoRow18
$ws.Range("F16").Value = $valorRow18
$ws.Range("E18").Value = $periodoRow16
$ws.Range("F18").Value = $valorRow16
